# Update the Amazon search-result data: replace each category's sample
# product title with a fresh search result, and record a result-count
# value of 2 in column C (error handling now reports how many results
# were found for each category search).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 1;  Text = "AmazonBasics Cat Activity Tree with Scratching Posts, Medium" },
    @{ Row = 2;  Text = "EZ Straw Just Straw Clean Processed Straw, Small Bale (1 cubic foot bale)" },
    @{ Row = 3;  Text = "AmazonBasics 92 Bright Multipurpose Copy Paper - 8.5 x 11 Inches, 10 Ream Case (5,000 Sheets)" },
    @{ Row = 4;  Text = "Perfect Cloud Double Airflow Memory Foam Pillow Featuring Cooling Ventilated Visco Foam Core and Mesh Trim for a Refreshing Sleeping Experience" },
    @{ Row = 5;  Text = "AmazonBasics 1/2-Inch Extra Thick Exercise Mat" },
    @{ Row = 6;  Text = "Dinner Forks,MCIRCO 12-Piece Good Quality Stainless Steel Table Forks Cutlery Set,8 Inch" },
    @{ Row = 7;  Text = "Utopia Kitchen 6 Pieces Bowl Set - Dishwasher Safe Opal Glassware - Microwave/Oven Friendly" },
    @{ Row = 8;  Text = "Amazon Brand - Solimo 18oz Disposable Plastic Party Cups, 200 Count, Red" },
    @{ Row = 9;  Text = "BEDSURE Sherpa Fleece Blanket Twin Size Red Plush Throw Blanket Fuzzy Soft Blanket Microfiber" },
    @{ Row = 10; Text = "No Item Found" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Text
    $ws.Cells.Item($r, 3).Value = 2
}
